$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "57.413.42"
$ws.Range("E2").Value = "  -0.05%  "

$ws.Range("D3").Value = "3.108.60"
$ws.Range("E3").Value = "  +0.06%  "

$ws.Range("D5").Value = "'525.54"
$ws.Range("E5").Value = "  +0.35%  "

$ws.Range("D6").Value = "'136.62"
$ws.Range("E6").Value = "  -3.22%  "

$ws.Range("E7").Value = "  +0.03%  "

$ws.Range("D8").Value = "3.109.53"
$ws.Range("E8").Value = "  +0.08%  "

$ws.Range("E9").Value = "  +2.56%  "

$ws.Range("D10").Value = "'7.20"
$ws.Range("E10").Value = "  -0.35%  "

$ws.Range("E11").Value = "  -0.84%  "

$ws.Range("D12").Value = "'0.396"
$ws.Range("E12").Value = "  +2.58%  "

$ws.Range("D13").Value = "3.647.94"
$ws.Range("E13").Value = "  +0.19%  "

$ws.Range("E14").Value = "  +2.73%  "

$ws.Range("D15").Value = "'25.28"
$ws.Range("E15").Value = "  -3.24%  "

$ws.Range("E16").Value = "  +0.24%  "

$ws.Range("D17").Value = "57.548.44"
$ws.Range("E17").Value = "  +0.01%  "

$ws.Range("D18").Value = "3.119.03"
$ws.Range("E18").Value = "  +0.38%  "

$ws.Range("D19").Value = "'5.94"
$ws.Range("E19").Value = "  -2.64%  "

$ws.Range("D20").Value = "'12.57"
$ws.Range("E20").Value = "  -1.73%  "

$ws.Range("D21").Value = "'7.91"
$ws.Range("E21").Value = "  -2.02%  "

$ws.Range("D22").Value = "'348.95"
$ws.Range("E22").Value = "  +3.76%  "

$ws.Range("E23").Value = "  -0.55%  "

$ws.Range("E24").Value = "  -0.09%  "

$ws.Range("D25").Value = "'68.27"
$ws.Range("E25").Value = "  +2.55%  "

$ws.Range("D26").Value = "'0.504"
$ws.Range("E26").Value = "  -1.76%  "

$ws.Range("D27").Value = "'0.167"
$ws.Range("E27").Value = "  -0.88%  "

$ws.Range("D28").Value = "'0.999"
$ws.Range("E28").Value = "  -0.24%  "

$ws.Range("D29").Value = "0.0₃0910"
$ws.Range("E29").Value = "  -0.92%  "

$ws.Range("D30").Value = "'7.41"
$ws.Range("E30").Value = "  +2.83%  "

$ws.Range("D31").Value = "'0.998"
$ws.Range("E31").Value = "  -0.07%  "

$ws.Range("E32").Value = "  +0.50%  "

$ws.Range("D33").Value = "'6.08"
$ws.Range("E33").Value = "  -6.83%  "

$ws.Range("D34").Value = "'21.09"
$ws.Range("E34").Value = "  +0.49%  "

$ws.Range("E35").Value = "  -2.63%  "

$ws.Range("D36").Value = "'4.94"
$ws.Range("E36").Value = "  +6.44%  "

$ws.Range("D37").Value = "'158.41"
$ws.Range("E37").Value = "  +0.54%  "

$ws.Range("D38").Value = "'6.15"
$ws.Range("E38").Value = "  +0.70%  "

$ws.Range("D39").Value = "'26.09"
$ws.Range("E39").Value = "  -3.18%  "

$ws.Range("D40").Value = "'1.24"
$ws.Range("E40").Value = "  -2.85%  "

$ws.Range("D41").Value = "'4.19"
$ws.Range("E41").Value = "  +6.16%  "

$ws.Range("D42").Value = "'0.0664"
$ws.Range("E42").Value = "  +0.46%  "

$ws.Range("E43").Value = "  +6.59%  "

$ws.Range("D44").Value = "'0.700"
$ws.Range("E44").Value = "  +1.95%  "

$ws.Range("D45").Value = "3.151.35"
$ws.Range("E45").Value = "  +0.12%  "

$ws.Range("E46").Value = "  -0.93%  "

$ws.Range("D47").Value = "2.343.10"
$ws.Range("E47").Value = "  +1.70%  "

$ws.Range("B48").Value = "FirstDigitalUSD"
$ws.Range("C48").Value = "https://coinranking.com/coin/cpjRxjFYD+firstdigitalusd-fdusd"
$ws.Range("D48").Value = "'1.00"
$ws.Range("E48").Value = "  +0.01%  "

$ws.Range("B49").Value = "VeChain"
$ws.Range("C49").Value = "https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet"
$ws.Range("D49").Value = "'0.0268"
$ws.Range("E49").Value = "  +3.27%  "

$ws.Range("D50").Value = "'0.956"
$ws.Range("E50").Value = "  -2.20%  "

$ws.Range("D51").Value = "'6.03"
$ws.Range("E51").Value = "  +0.16%  "

